$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '37.694.83'
$ws.Range("E2").Value = '  -0.02%  '

# Row 3
$ws.Range("D3").Value = '2.083.38'
$ws.Range("E3").Value = '  +0.43%  '

# Row 4
$ws.Range("E4").Value = '  +0.02%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '233.32'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.06%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.637'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.32%  '

# Row 7
$ws.Range("E7").Value = '  +0.02%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '57.93'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.25%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.391'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.17%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0778'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.50%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.109'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +2.99%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.23%  '

# Row 13
$ws.Range("D13").Value = '2.390.77'
$ws.Range("E13").Value = '  +0.45%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.02'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.94%  '

# Row 15
$ws.Range("E15").Value = '  +0.30%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.34'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.89%  '

# Row 17
$ws.Range("D17").Value = '2.101.15'
$ws.Range("E17").Value = '  +1.27%  '

# Row 18
$ws.Range("D18").Value = '37.694.60'
$ws.Range("E18").Value = '  +0.12%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.07'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.34%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.85'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.24%  '

# Row 21
$ws.Range("E21").Value = '  +0.06%  '

# Row 22
$ws.Range("E22").Value = '  +0.42%  '

# Row 23
$ws.Range("E23").Value = '  -0.07%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -1.48%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.38'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.11%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.66'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.43%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '170.30'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.52%  '

# Row 28
$ws.Range("E28").Value = '  -4.19%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.50'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.42%  '

# Row 30
$ws.Range("E30").Value = '  -0.25%  '

# Row 31
$ws.Range("E31").Value = '  +1.46%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.22%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0635'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.67%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.50%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.49'
$ws.Range("D35").Style = "Normal"

# Row 36
$ws.Range("E36").Value = '  -0.44%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.31'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.26%  '

# Row 38
$ws.Range("E38").Value = '  +0.02%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.35'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.01%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0236'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.68%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '100.97'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.15%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0958'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.46%  '

# Row 43
$ws.Range("B43").Value = 'TrustWalletToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.20'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.17%  '

# Row 44
$ws.Range("B44").Value = 'HuobiToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.73%  '

# Row 45
$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '1.462.27'
$ws.Range("E45").Value = '  +0.46%  '

# Row 46
$ws.Range("B46").Value = 'InjectiveProtocol'
$ws.Range("C46").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.65'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.12%  '

# Row 47
$ws.Range("E47").Value = '  -0.93%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.04'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.54%  '

# Row 49
$ws.Range("E49").Value = '  -2.31%  '

# Row 50
$ws.Range("E50").Value = '  -1.85%  '

# Row 51
$ws.Range("D51").Value = '2.275.80'
$ws.Range("E51").Value = '  +0.46%  '
